$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns and fix capitalization of "de/del/la/los" in place names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B3").Value = "Rincón De Romos"
$ws.Range("B12").Value = "Hidalgo Del Parral"
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("A21").Value = "Estado De México"
$ws.Range("B22").Value = "Tlalnepantla De Baz"
$ws.Range("B26").Value = "Silao De La Victoria"
$ws.Range("B28").Value = "Cuautepec De Hinojosa"
$ws.Range("B30").Value = "Pachuca De Soto"
$ws.Range("B32").Value = "Tulancingo De Bravo"
$ws.Range("B36").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B40").Value = "Ojuelos De Jalisco"
$ws.Range("B43").Value = "Tepatitlán De Morelos"
$ws.Range("B47").Value = "Valle De Guadalupe"
$ws.Range("B61").Value = "Bahía De Banderas"
$ws.Range("B67").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B68").Value = "Oaxaca De Juárez"
$ws.Range("B70").Value = "Tlacolula De Matamoros"
$ws.Range("B82").Value = "Nacozari De García"
$ws.Range("B84").Value = "Amatlán De Los Reyes"
$ws.Range("B85").Value = "Ignacio De La Llave"

# Remove the trailing metadata/footer rows (95-99) that are no longer part of the clean data range
$ws.Rows("95:99").Delete()

